# TimeLog_EricCarlson.xlsx -- "Week 8" sheet update
# The 3D World is fixed! Adds the missing "Type" entries for a couple of
# existing rows, fills in the rest of Saturday's and Sunday's work log
# (which were previously blank template rows), and adds the week's
# grand-total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 8")

# ---------------------------------------------------------------------
# 1) Friday section (rows 19-22): the "Type" column (C) was left blank
#    for the single logged entry -- fill it in.
# ---------------------------------------------------------------------
$ws.Range("C21").Value2 = "Fixing 3D World"

# ---------------------------------------------------------------------
# 2) Saturday section (rows 25-30): the "Type" column was blank for the
#    first two entries, and a third entry plus the "Today's Total" row
#    still need to be added.
# ---------------------------------------------------------------------
$ws.Range("C27").Value2 = "Fixing 3D World"
$ws.Range("C28").Value2 = "Fixing 3D World"

# New data row 29 -- copy formatting from the row above (row 28).
$ws.Range("A28:D28").Copy($ws.Range("A29:D29"))
$ws.Range("A29").Value2 = 0.54166666666666663
$ws.Range("B29").Value2 = 0.625
$ws.Range("C29").Value2 = "Fixing 3D World"
$ws.Range("D29").Value2 = 2

# New "Today's Total" row 30 -- copy formatting from the equivalent row
# in the Friday section above (row 22), which carries the merged /
# centred style.
$ws.Range("A22:D22").Copy($ws.Range("A30:D30"))
$ws.Range("D30").Value2 = 3.5

# ---------------------------------------------------------------------
# 3) Sunday header block originally sat at rows 38-39; make room for the
#    new Saturday data above it by inserting a single row before it.
#    (Clear the clipboard/marching-ants state left by the Copy() calls
#    above first, otherwise a bare Insert() pastes the clipboard
#    contents instead of just shifting cells down.)
# ---------------------------------------------------------------------
$excel.CutCopyMode = $false
$ws.Rows.Item(38).Insert()

# ---------------------------------------------------------------------
# 4) Today (Saturday) section continues -- now a new empty Week 9 data
#    block exists at rows 34-37 (the old blank rows 34-37 before the
#    insert). Populate it.
# ---------------------------------------------------------------------
$ws.Range("A28:D28").Copy($ws.Range("A34:D34"))
$ws.Range("A34").Value2 = 0.5
$ws.Range("B34").Value2 = 0.5625
$ws.Range("C34").Value2 = "RC Car Bug Fixing"
$ws.Range("D34").Value2 = 1.5

$ws.Range("A28:D28").Copy($ws.Range("A35:D35"))
$ws.Range("A35").Value2 = 0.58333333333333337
$ws.Range("B35").Value2 = 0.66666666666666663
$ws.Range("C35").Value2 = "Practicing Presentation"
$ws.Range("D35").Value2 = 2

$ws.Range("A28:D28").Copy($ws.Range("A36:D36"))
$ws.Range("A36").Value2 = 0.72916666666666663
$ws.Range("B36").Value2 = 0.75
$ws.Range("C36").Value2 = "Fixing 3D World"
$ws.Range("D36").Value2 = 0.5

$ws.Range("A22:D22").Copy($ws.Range("A37:D37"))
$ws.Range("D37").Value2 = 4

# ---------------------------------------------------------------------
# 5) Sunday section now lives at rows 39 (header), 40 (table header),
#    41-42 (data, newly blank rows) and 43 (Today's Total, newly blank).
# ---------------------------------------------------------------------
$ws.Range("A28:D28").Copy($ws.Range("A41:D41"))
$ws.Range("A41").Value2 = 0.91666666666666663
$ws.Range("B41").Value2 = 0.95833333333333337
$ws.Range("C41").Value2 = "Logo/Presentation"
$ws.Range("D41").Value2 = 1

$ws.Range("A28:D28").Copy($ws.Range("A42:D42"))
$ws.Range("A42").Value2 = 0.95833333333333337
$ws.Range("B42").Value2 = 0.16666666666666666
$ws.Range("C42").Value2 = "Fixing 3D World: Fixed!"
$ws.Range("D42").Value2 = 5

$ws.Range("A22:D22").Copy($ws.Range("A43:D43"))
$ws.Range("D43").Value2 = 6

# ---------------------------------------------------------------------
# 6) New "Week's Total" row 45 (row 44 left blank, matching the blank
#    separator row used before every other section header).
# ---------------------------------------------------------------------
$ws.Range("A22:D22").Copy($ws.Range("A45:D45"))
$ws.Range("A45").Value2 = "Week's Total"
$ws.Range("D45").Formula = "=SUM(D10,D17,D22,D30,D37,D43)"

# ---------------------------------------------------------------------
# 7) Column C needs to widen to fit the new, longer "Type" entries.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19

# ---------------------------------------------------------------------
# 8) Update the view state to match where the user ended up editing.
# ---------------------------------------------------------------------
$ws.Range("D45").Select()
